$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row: "q" -> "question", "a" -> "answer"
$ws.Range("A1").Value = "question"
$ws.Range("B1").Value = "answer"

# Update selection to B1
$ws.Range("B1").Select()
